# Bug-Report.xlsx: mark both logged bugs as "Closed" in the Status column (K)
# and leave the selection where the author last clicked (N2), matching the
# commit "Doc: updated Bugs status to closed".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "Closed"
$ws.Range("K3").Value = "Closed"

$ws.Range("N2").Select() | Out-Null
